$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# dSF (column F) corrections from repulled data
$ws.Range("F6").Value = -6
$ws.Range("F9").Value = -4
$ws.Range("F15").Value = -1
$ws.Range("F18").Value = 7
$ws.Range("F19").Value = -7
$ws.Range("F27").Value = -3
$ws.Range("F31").Value = -1
$ws.Range("F33").Value = 3
$ws.Range("F35").Value = -1
$ws.Range("F38").Value = -2
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = -1
$ws.Range("F43").Value = 2
$ws.Range("F46").Value = -1
$ws.Range("F48").Value = 1
$ws.Range("F52").Value = -1
$ws.Range("F54").Value = -7
$ws.Range("F55").Value = -3
$ws.Range("F68").Value = -3

# dS0 (column E) correction
$ws.Range("E56").Value = 2
